# Auto-generated Excel COM-interop edit script
# Applies cell-value corrections to the "Sagittarius_Profits" workbook sheets
# (market-price / leve-profit recompute), per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("I40").Value2 = 1702.871
$ws.Range("H40").Value2 = 1894.907
$ws.Range("M40").Value2 = -1527.871
$ws.Range("J40").Value2 = 2391
$ws.Range("L40").Value2 = 2391
$ws.Range("K40").Value2 = 1702.871
$ws.Range("N40").Value2 = -2741
# row 86
$ws.Range("K86").Value2 = 1483.1666
$ws.Range("M86").Value2 = -360.1666
$ws.Range("I86").Value2 = 1483.1666
$ws.Range("H86").Value2 = 1833
# row 89
$ws.Range("H89").Value2 = 1833
$ws.Range("M89").Value2 = -1799.833000000001
$ws.Range("K89").Value2 = 7415.833000000001
$ws.Range("I89").Value2 = 1483.1666
# row 113
$ws.Range("N113").Value2 = -10385.7778
$ws.Range("H113").Value2 = 3940
$ws.Range("J113").Value2 = 3877.7778
$ws.Range("L113").Value2 = 3877.7778
# row 132
$ws.Range("K132").Value2 = 7112.1819
$ws.Range("M132").Value2 = -4582.1819
$ws.Range("I132").Value2 = 2370.7273
$ws.Range("H132").Value2 = 2467.4614
# row 135
$ws.Range("J135").Value2 = 7073.6665
$ws.Range("N135").Value2 = -68732.9985
$ws.Range("L135").Value2 = 63662.9985
$ws.Range("K135").Value2 = 36469.8
$ws.Range("M135").Value2 = -33934.8
$ws.Range("I135").Value2 = 4052.2
$ws.Range("H135").Value2 = 4749.4614
# row 137
$ws.Range("H137").Value2 = 4407.3076
$ws.Range("K137").Value2 = 6932.3334
$ws.Range("M137").Value2 = -4382.3334
$ws.Range("I137").Value2 = 2310.7778

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("K61").Value2 = 2963.889
$ws.Range("M61").Value2 = -2751.889
$ws.Range("I61").Value2 = 2963.889
$ws.Range("H61").Value2 = 3167.5
# row 74
$ws.Range("K74").Value2 = 1374.9
$ws.Range("M74").Value2 = -500.9000000000001
$ws.Range("I74").Value2 = 1374.9
$ws.Range("H74").Value2 = 1409.0714
# row 77
$ws.Range("M77").Value2 = -2506.5
$ws.Range("K77").Value2 = 6874.5
$ws.Range("I77").Value2 = 1374.9
$ws.Range("H77").Value2 = 1409.0714
# row 136
$ws.Range("K136").Value2 = 8891.667000000001
$ws.Range("M136").Value2 = -6341.667000000001
$ws.Range("I136").Value2 = 2963.889
$ws.Range("H136").Value2 = 3167.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("K22").Value2 = 223.14285
$ws.Range("M22").Value2 = -50.14285000000001
$ws.Range("I22").Value2 = 223.14285
$ws.Range("H22").Value2 = 220.5
# row 86
$ws.Range("J86").Value2 = 2296.3333
$ws.Range("L86").Value2 = 2296.3333
$ws.Range("N86").Value2 = -4542.3333
$ws.Range("K86").Value2 = 2003
$ws.Range("M86").Value2 = -880
$ws.Range("I86").Value2 = 2003
$ws.Range("H86").Value2 = 2179
# row 89
$ws.Range("H89").Value2 = 2179
$ws.Range("J89").Value2 = 2296.3333
$ws.Range("L89").Value2 = 11481.6665
$ws.Range("M89").Value2 = -4399
$ws.Range("N89").Value2 = -22713.6665
$ws.Range("K89").Value2 = 10015
$ws.Range("I89").Value2 = 2003
# row 94
$ws.Range("M94").Value2 = -1265.3334
$ws.Range("I94").Value2 = 1716.3334
$ws.Range("H94").Value2 = 1716.3334
$ws.Range("K94").Value2 = 1716.3334
# row 97
$ws.Range("L97").Value2 = 19900
$ws.Range("K97").Value2 = 9977
$ws.Range("M97").Value2 = -8986
$ws.Range("N97").Value2 = -21882
$ws.Range("I97").Value2 = 9977
$ws.Range("H97").Value2 = 14938.5
$ws.Range("J97").Value2 = 19900
# row 99
$ws.Range("K99").Value2 = 999.5
$ws.Range("M99").Value2 = 498.5
$ws.Range("I99").Value2 = 999.5
$ws.Range("H99").Value2 = 1066
# row 107
$ws.Range("K107").Value2 = 1384
$ws.Range("M107").Value2 = 536
$ws.Range("I107").Value2 = 1384
$ws.Range("H107").Value2 = 1435.5454
# row 134
$ws.Range("H134").Value2 = 1840.3334
$ws.Range("M134").Value2 = -2986.0002
$ws.Range("K134").Value2 = 5521.0002
$ws.Range("I134").Value2 = 1840.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("K58").Value2 = 1786.5555
$ws.Range("M58").Value2 = -1583.5555
$ws.Range("I58").Value2 = 1786.5555
$ws.Range("H58").Value2 = 1929.4546
# row 122
$ws.Range("M122").Value2 = 333.5
$ws.Range("K122").Value2 = 2116.5
$ws.Range("I122").Value2 = 705.5
$ws.Range("H122").Value2 = 1263.25
# row 132
$ws.Range("K132").Value2 = 8268.900000000001
$ws.Range("M132").Value2 = -5738.900000000001
$ws.Range("I132").Value2 = 2756.3
$ws.Range("H132").Value2 = 2755
# row 134
$ws.Range("H134").Value2 = 1809.1666
$ws.Range("J134").Value2 = 4000
$ws.Range("L134").Value2 = 12000
$ws.Range("M134").Value2 = -2295
$ws.Range("N134").Value2 = -17070
$ws.Range("K134").Value2 = 4830
$ws.Range("I134").Value2 = 1610
# row 136
$ws.Range("K136").Value2 = 5359.666499999999
$ws.Range("M136").Value2 = -2809.666499999999
$ws.Range("I136").Value2 = 1786.5555
$ws.Range("H136").Value2 = 1929.4546
# row 141
$ws.Range("J141").Value2 = 38036.855
$ws.Range("L141").Value2 = 38036.855
$ws.Range("N141").Value2 = -48396.855
$ws.Range("M141").ClearContents()
$ws.Range("K141").Value2 = 0
$ws.Range("I141").Value2 = 0
$ws.Range("H141").Value2 = 38036.855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 23
$ws.Range("K23").Value2 = 300000240
$ws.Range("M23").Value2 = -300000005
$ws.Range("N23").Value2 = -2196.8
$ws.Range("I23").Value2 = 100000080
$ws.Range("H23").Value2 = 37500388
$ws.Range("J23").Value2 = 575.6
$ws.Range("L23").Value2 = 1726.8
# row 137
$ws.Range("H137").Value2 = 2043.1333
$ws.Range("J137").Value2 = 2368.4546
$ws.Range("L137").Value2 = 7105.3638
$ws.Range("N137").Value2 = -17305.3638

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 10
$ws.Range("K10").Value2 = 1000000
$ws.Range("M10").Value2 = -999831
$ws.Range("I10").Value2 = 1000000
$ws.Range("H10").Value2 = 503400
# row 15
$ws.Range("H15").Value2 = 44987
$ws.Range("J15").Value2 = 44987
$ws.Range("L15").Value2 = 44987
$ws.Range("N15").Value2 = -45563
# row 70
$ws.Range("J70").Value2 = 5000
$ws.Range("L70").Value2 = 5000
$ws.Range("M70").Value2 = -4041.25
$ws.Range("N70").Value2 = -5540
$ws.Range("K70").Value2 = 4311.25
$ws.Range("I70").Value2 = 4311.25
$ws.Range("H70").Value2 = 4449
# row 73
$ws.Range("H73").Value2 = 4449
$ws.Range("J73").Value2 = 5000
$ws.Range("L73").Value2 = 5000
$ws.Range("K73").Value2 = 4311.25
$ws.Range("M73").Value2 = -3375.25
$ws.Range("N73").Value2 = -6872
$ws.Range("I73").Value2 = 4311.25
# row 74
$ws.Range("L74").Value2 = 50000
$ws.Range("N74").Value2 = -51872
$ws.Range("H74").Value2 = 50000
$ws.Range("J74").Value2 = 50000
# row 77
$ws.Range("L77").Value2 = 150000
$ws.Range("N77").Value2 = -159360
$ws.Range("H77").Value2 = 50000
$ws.Range("J77").Value2 = 50000
# row 80
$ws.Range("K80").Value2 = 2999.5
$ws.Range("M80").Value2 = -2001.5
$ws.Range("I80").Value2 = 2999.5
$ws.Range("H80").Value2 = 2999.5
# row 81
$ws.Range("N81").Value2 = -46983
$ws.Range("H81").Value2 = 44987
$ws.Range("J81").Value2 = 44987
$ws.Range("L81").Value2 = 44987
# row 83
$ws.Range("K83").Value2 = 14997.5
$ws.Range("M83").Value2 = -10005.5
$ws.Range("I83").Value2 = 2999.5
$ws.Range("H83").Value2 = 2999.5
# row 84
$ws.Range("N84").Value2 = -144945
$ws.Range("H84").Value2 = 44987
$ws.Range("J84").Value2 = 44987
$ws.Range("L84").Value2 = 134961
# row 132
$ws.Range("K132").Value2 = 18868.221
$ws.Range("M132").Value2 = -16338.221
$ws.Range("I132").Value2 = 6289.407
$ws.Range("H132").Value2 = 5946.9395
# row 134
$ws.Range("H134").Value2 = 32264.8
$ws.Range("J134").Value2 = 32264.8
$ws.Range("L134").Value2 = 96794.39999999999
$ws.Range("N134").Value2 = -101864.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("L7").Value2 = 8890.182000000001
$ws.Range("K7").Value2 = 4004
$ws.Range("M7").Value2 = -3892
$ws.Range("N7").Value2 = -9114.182000000001
$ws.Range("I7").Value2 = 4004
$ws.Range("H7").Value2 = 8483
$ws.Range("J7").Value2 = 8890.182000000001
# row 40
$ws.Range("I40").Value2 = 3386.25
$ws.Range("H40").Value2 = 4192
$ws.Range("M40").Value2 = -3250.25
$ws.Range("J40").Value2 = 5266.3335
$ws.Range("L40").Value2 = 5266.3335
$ws.Range("K40").Value2 = 3386.25
$ws.Range("N40").Value2 = -5538.3335
# row 46
$ws.Range("K46").Value2 = 1173.7059
$ws.Range("M46").Value2 = -985.7058999999999
$ws.Range("I46").Value2 = 1173.7059
$ws.Range("H46").Value2 = 1390.8823
# row 68
$ws.Range("K68").Value2 = 1199.2
$ws.Range("M68").Value2 = -450.2
$ws.Range("I68").Value2 = 1199.2
$ws.Range("H68").Value2 = 1199.2
# row 71
$ws.Range("K71").Value2 = 5996
$ws.Range("M71").Value2 = -2252
$ws.Range("I71").Value2 = 1199.2
$ws.Range("H71").Value2 = 1199.2
# row 122
$ws.Range("L122").Value2 = 22384.125
$ws.Range("M122").Value2 = -9534.0772
$ws.Range("N122").Value2 = -27284.125
$ws.Range("K122").Value2 = 11984.0772
$ws.Range("I122").Value2 = 3994.6924
$ws.Range("H122").Value2 = 5907.3447
$ws.Range("J122").Value2 = 7461.375
# row 126
$ws.Range("N126").Value2 = -31610.546
$ws.Range("K126").Value2 = 12012
$ws.Range("M126").Value2 = -9542
$ws.Range("I126").Value2 = 4004
$ws.Range("H126").Value2 = 8483
$ws.Range("J126").Value2 = 8890.182000000001
$ws.Range("L126").Value2 = 26670.546
# row 136
$ws.Range("K136").Value2 = 6299.25
$ws.Range("M136").Value2 = -3749.25
$ws.Range("I136").Value2 = 2099.75
$ws.Range("H136").Value2 = 2099.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 68
$ws.Range("N68").Value2 = -21619
$ws.Range("H68").Value2 = 21249.25
$ws.Range("J68").Value2 = 19997
$ws.Range("L68").Value2 = 19997
# row 71
$ws.Range("N71").Value2 = -68103
$ws.Range("L71").Value2 = 59991
$ws.Range("H71").Value2 = 21249.25
$ws.Range("J71").Value2 = 19997
# row 136
$ws.Range("K136").Value2 = 7578.999899999999
$ws.Range("M136").Value2 = -5028.999899999999
$ws.Range("I136").Value2 = 2526.3333
$ws.Range("H136").Value2 = 2879.7144
